$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row - Right column (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row - Right column (B12): 42 -> 70
$ws.Range("B12").Value = 70

# Update "Total" row - Max column (E12): "42/84" -> "70/140"
$ws.Range("E12").Value = "70/140"
